# Increase line spacing in the SKILLS bullet list by adding
# "space after" (w:spacing w:after="40") = 2 pt to each list item
# that already carries the explicit w:ind left="454" hanging="170"
# formatting (i.e. every bullet except the first, empty placeholder
# bullet immediately under the heading).

$d = $word.ActiveDocument

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)

    if ($p.Range.ListFormat.ListType -ne 0 -and [Math]::Round($p.Format.LeftIndent) -eq 23) {
        $p.Format.SpaceAfter = 2
    }
}
